$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 2 new rows at row 616, pushing existing rows 616-649 down to 618-651
$ws.Range("A616:A617").EntireRow.Insert()

# Fill the new row 616
$ws.Range("A616").Value = 3
$ws.Range("B616").Value = "Femacal de La Calera"
$ws.Range("C616").Value = "Coquimbo"
$ws.Range("D616").Value = 44516
$ws.Range("E616").Value = 5
$ws.Range("F616").Value = "Fruta"
$ws.Range("G616").Value = 100102
$ws.Range("H616").Value = "Cítricos"
$ws.Range("I616").Value = 100102003
$ws.Range("J616").Value = "Limón"
$ws.Range("K616").Value = "Sin especificar"
$ws.Range("L616").Value = "1a amarillo"
$ws.Range("M616").Value = 398
$ws.Range("N616").Value = 4500
$ws.Range("O616").Value = 6000
$ws.Range("P616").Value = 5254
$ws.Range("Q616").Value = "$/malla 16 kilos"
$ws.Range("R616").Value = "Provincia de Quillota"
$ws.Range("S616").Value = 328
$ws.Range("T616").Value = 16

# Fill the new row 617
$ws.Range("A617").Value = 3
$ws.Range("B617").Value = "Femacal de La Calera"
$ws.Range("C617").Value = "Coquimbo"
$ws.Range("D617").Value = 44516
$ws.Range("E617").Value = 5
$ws.Range("F617").Value = "Fruta"
$ws.Range("G617").Value = 100102
$ws.Range("H617").Value = "Cítricos"
$ws.Range("I617").Value = 100102003
$ws.Range("J617").Value = "Limón"
$ws.Range("K617").Value = "Sin especificar"
$ws.Range("L617").Value = "2a amarillo"
$ws.Range("M617").Value = 275
$ws.Range("N617").Value = 3500
$ws.Range("O617").Value = 4000
$ws.Range("P617").Value = 3727
$ws.Range("Q617").Value = "$/malla 16 kilos"
$ws.Range("R617").Value = "Provincia de Quillota"
$ws.Range("S617").Value = 233
$ws.Range("T617").Value = 16

# Ensure D616/D617 carry the same date number format as other D cells (style s="2")
$ws.Range("D618").Copy()
$ws.Range("D616:D617").PasteSpecial(-4122) # xlPasteFormats
